# Apply cryptos.xlsx price/volume/name updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.383.67'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '1.716.60'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '225.37'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5305'
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.06705'
$ws.Range('E8').Value = '  +1.93%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2667'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.94'
$ws.Range('E10').Value = '  -2.89%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07705'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.510'
$ws.Range('E12').Value = '  -2.03%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.951.26'
$ws.Range('E13').Value = '  -0.63%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.714.33'
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5865'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '0.0₅8230'
$ws.Range('E16').Value = '  -0.58%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '68.05'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '27.358.22'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '223.49'
$ws.Range('E19').Value = '  +2.61%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.667'
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.52'
$ws.Range('E22').Value = '  -0.73%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.044'
$ws.Range('E23').Value = '  -0.48%  '
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '144.59'
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.691'
$ws.Range('E26').Value = '  -3.59%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1210'
$ws.Range('E27').Value = '  -2.00%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.250'
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '16.29'
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05367'
$ws.Range('E30').Value = '  -2.18%  '
$ws.Range('E31').Value = '  -0.59%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.483'
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.636'
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.872'
$ws.Range('E35').Value = '  +0.32%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.9586'
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.389'
$ws.Range('E37').Value = '  -1.53%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5886'
$ws.Range('E38').Value = '  -1.31%  '
$ws.Range('D39').Value = '1.149.55'
$ws.Range('E39').Value = '  +9.09%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01647'
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.799'
$ws.Range('E41').Value = '  -1.70%  '
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8422'
$ws.Range('E43').Value = '  -1.19%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '100.89'
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('D45').Value = '1.857.69'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('E46').Value = '  -3.13%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '57.74'
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.4586'
$ws.Range('E48').Value = '  +2.47%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.005'
$ws.Range('E49').Value = '  +0.18%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.114'
$ws.Range('E50').Value = '  -1.14%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05196'
$ws.Range('E51').Value = '  -0.86%  '
